$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C ("Förändrad") for rows 2-18 from 45204 to 45207
for ($row = 2; $row -le 18; $row++) {
    $ws.Cells.Item($row, 3).Value = 45207
}
